# Add files via upload
# Adds 5 new scenario columns (AND/incentive, Nested-IF outstanding, OR/warning,
# AND/manager, Nested-IF bonus %) to the "Scanario based que on CONDITIONS" sheet,
# plus their matching question rows (11-15) further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 8): new headers in S8:W8, same look as the existing
#    header cells (copy format from R8, which already carries the header style).
# ---------------------------------------------------------------------------
$ws.Rows("8").RowHeight = 189

$ws.Range("S8").Value = 'AND=if(and(exp>4,sals>9500),"incentive","no incentive"'
$ws.Range("T8").Value = 'Nested if simple = rating>=4.5,"outstanding",rating>4,"very good",rating>3,"satisfactory",rating<3,"poor"'
$ws.Range("U8").Value = 'OR = sales>9000,exp<2,"warning"'
$ws.Range("V8").Value = 'AND = exp>=5,dept<>"HR" (not sign=<>)'
$ws.Range("W8").Value = 'Nested if simple= sales>=12000,"15%",sales>=10000,"10%",sales>=9000,"5%",sales<9000,"0"'

$ws.Range("R8").Copy()
$ws.Range("S8:W8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Data / formula rows 10-15 for the five new columns. Formulas (and the
#    occasional "typo" left by the original author, e.g. U11 referencing row
#    10) are reproduced exactly as authored.
# ---------------------------------------------------------------------------

# Column S - IF(AND(...)) incentive check
$ws.Range("S10").Formula = '=IF(AND(E10>4,C10>9500),"incentive","no incentive")'
$ws.Range("S11").Formula = '=IF(AND(E11>4,C11>=9500),"incentive","no incentive")'
$ws.Range("S12").Formula = '=IF(AND(E12>4,C12>9500),"incentive","no incentive")'
$ws.Range("S13").Formula = '=IF(AND(E13>4,C13>9500),"incentive","no incentive")'
$ws.Range("S14").Formula = '=IF(AND(E14>4,C14>9500),"incentive","no incentive")'
$ws.Range("S15").Formula = '=IF(AND(E15>4,C15>9500),"incentive","no incentive")'

# Column T - nested IFS rating scale (array-entered, as in the original file)
$ws.Range("T10").FormulaArray = '=IFS(F10>=4.5,"outstanding",F10>=4,"very good",F10>=3,"satisfactory",F10<3,"poor")'
$ws.Range("T11").FormulaArray = '=IFS(F11>=4.5,"outstanding",F11>=4,"very good",F11>=3,"satisfactory",F11<3,"poor")'
$ws.Range("T12").FormulaArray = '=IFS(F12>=4.5,"outstanding",F12>=4,"very good",F12>=3,"satisfactory",F12<3,"poor")'
$ws.Range("T13").FormulaArray = '=IFS(F13>=4.5,"outstanding",F13>=4,"very good",F13>=3,"satisfactory",F13<3,"poor")'
$ws.Range("T14").FormulaArray = '=IFS(F14>=4.5,"outstanding",F14>=4,"very good",F14>=3,"satisfactory",F14<3,"poor")'
$ws.Range("T15").FormulaArray = '=IFS(F15>=4.5,"outstanding",F15>=4,"very good",F15>=3,"satisfactory",F15<3,"poor")'

# Column U - IF(OR(...)) sales/experience warning
$ws.Range("U10").Formula = '=IF(OR(C10<9000,E10<2),"warning","NO")'
$ws.Range("U11").Formula = '=IF(OR(C10<9000,E10<2),"WARNING","NO")'
$ws.Range("U12").Formula = '=IF(OR(C12<9000,E12<2),"WARNING","NO")'
$ws.Range("U13").Formula = '=IF(OR(C13<9000,E13<2),"WARNING","NO")'
$ws.Range("U14").Formula = '=IF(OR(C14<9000,E14<2),"WARNING","NO")'
$ws.Range("U15").Formula = '=IF(OR(C15<9000,E15<2),"WARNING","NO")'

# Column V - IF(AND(...)) experience/department manager check
$ws.Range("V10").Formula = '=IF(AND(E10>=5,B10<>"HR"),"manager","NO")'
$ws.Range("V11").Formula = '=IF(AND(E11>=5,B11<>"HR"),"manager","no")'
$ws.Range("V12").Formula = '=IF(AND(E12>=5,B12<>"HR"),"manager","no")'
$ws.Range("V13").Formula = '=IF(AND(E13>=5,B13<>"HR"),"manager","no")'
$ws.Range("V14").Formula = '=IF(AND(E14>=5,B14<>"HR"),"manager","no")'
$ws.Range("V15").Formula = '=IF(AND(E15>=5,B15<>"HR"),"manager","no")'

# Column W - nested IFS bonus percentage (array-entered)
$ws.Range("W10").FormulaArray = '=IFS(C10>=12000,"15%",C10>=10000,"10%",C10>=9000,"5%",C10<9000,"0")'
$ws.Range("W11").FormulaArray = '=IFS(C11>=12000,"15%",C11>=10000,"10%",C11>=9000,"5%",C11<9000,"0")'
$ws.Range("W12").FormulaArray = '=IFS(C12>=12000,"15%",C12>=10000,"10%",C12>=9000,"5%",C12<9000,"0")'
$ws.Range("W13").FormulaArray = '=IFS(C13>=12000,"15%",C13>=10000,"10%",C13>=9000,"5%",C13<9000,"0")'
$ws.Range("W14").FormulaArray = '=IFS(C14>=12000,"15%",C14>=10000,"10%",C14>=9000,"5%",C14<9000,"0")'
$ws.Range("W15").FormulaArray = '=IFS(C15>=12000,"15%",C15>=10000,"10%",C15>=9000,"5%",C15<9000,"0")'

# ---------------------------------------------------------------------------
# 3. New question rows beneath the existing Q&A list (rows 17..35), matching
#    the sheet's "only even data rows have text" pattern.
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = '11)show "incentive" if the emp has more than 4 yr exp and sales are over $9500.otherwise show "no incentive"'
$ws.Range("A39").Value = '12)4.5 and above -"outstanding",4.0-4.9="very good",3.0-3.99="satisfactory",blow 3.0 ="poor"'
$ws.Range("A41").Value = '13)if sales are less than$9000 or exp is under 2 yrs ,"warning"'
$ws.Range("A43").Value = '14)if exp is 5 or more than that and dept is not HR ,show "manager"'
$ws.Range("A45").Value = '15) assign bonus percentage based on sales,>=12000="15%",>=10000="10%",>=9000="5%",below 9000="0%"'

# ---------------------------------------------------------------------------
# 4. Column widths for the new columns (best-effort match of the authored
#    widths) and refreshed view/selection state.
# ---------------------------------------------------------------------------
$ws.Range("S1").ColumnWidth = 11.74
$ws.Range("T1").ColumnWidth = 23.74
$ws.Range("U1").ColumnWidth = 11.31
$ws.Range("V1").ColumnWidth = 8.17
$ws.Range("W1").ColumnWidth = 14.88

$ws.Range("M1").Select()
